# Update cryptocurrency price/volume data in the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.178.48"
$ws.Range("E2").Value = "  +0.79%  "

$ws.Range("D3").Value = "'3.568.31"
$ws.Range("E3").Value = "  -0.80%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").Value = "'609.52"
$ws.Range("E5").Value = "  -0.04%  "

$ws.Range("D6").Value = "'146.53"
$ws.Range("E6").Value = "  -1.15%  "

$ws.Range("D7").Value = "'3.570.62"
$ws.Range("E7").Value = "  -0.71%  "

$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.16%  "

$ws.Range("D9").Value = "'0.514"
$ws.Range("E9").Value = "  +5.17%  "

$ws.Range("D10").Value = "'7.90"
$ws.Range("E10").Value = "  -1.96%  "

$ws.Range("E11").Value = "  -2.20%  "

$ws.Range("D12").Value = "'0.414"
$ws.Range("E12").Value = "  -0.18%  "

$ws.Range("D13").Value = "'4.168.12"
$ws.Range("E13").Value = "  -0.90%  "

$ws.Range("D14").Value = "'0.0000197"
$ws.Range("E14").Value = "  -5.90%  "

$ws.Range("D15").Value = "'29.17"
$ws.Range("E15").Value = "  -2.65%  "

$ws.Range("D16").Value = "'3.562.14"
$ws.Range("E16").Value = "  -2.22%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "'66.912.47"
$ws.Range("E17").Value = "  +0.29%  "

$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "'0.118"
$ws.Range("E18").Value = "  +1.25%  "

$ws.Range("D19").Value = "'11.13"
$ws.Range("E19").Value = "  -3.28%  "

$ws.Range("D20").Value = "'6.24"
$ws.Range("E20").Value = "  -1.50%  "

$ws.Range("D21").Value = "'14.74"
$ws.Range("E21").Value = "  -2.40%  "

$ws.Range("D22").Value = "'428.80"
$ws.Range("E22").Value = "  +0.22%  "

$ws.Range("D23").Value = "'0.601"
$ws.Range("E23").Value = "  -3.16%  "

$ws.Range("D24").Value = "'77.85"
$ws.Range("E24").Value = "  -1.39%  "

$ws.Range("D25").Value = "'3.710.99"
$ws.Range("E25").Value = "  -0.62%  "

$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").Value = "'0.0000118"
$ws.Range("E27").Value = "  -3.68%  "

$ws.Range("D28").Value = "'8.11"
$ws.Range("E28").Value = "  -2.41%  "

$ws.Range("E29").Value = "  -1.06%  "

$ws.Range("D30").Value = "'9.10"
$ws.Range("E30").Value = "  -2.44%  "

$ws.Range("D31").Value = "'0.992"
$ws.Range("E31").Value = "  -0.80%  "

$ws.Range("D32").Value = "'3.575.31"
$ws.Range("E32").Value = "  -0.53%  "

$ws.Range("E33").Value = "  -1.40%  "

$ws.Range("D34").Value = "'24.55"
$ws.Range("E34").Value = "  -3.70%  "

$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").Value = "'1.37"
$ws.Range("E35").Value = "  -6.56%  "

$ws.Range("B36").Value = "USDe"
$ws.Range("C36").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.00%  "

$ws.Range("D37").Value = "'7.73"
$ws.Range("E37").Value = "  -1.53%  "

$ws.Range("E38").Value = "  -3.09%  "

$ws.Range("D39").Value = "'177.79"
$ws.Range("E39").Value = "  +0.13%  "

$ws.Range("D40").Value = "'5.34"
$ws.Range("E40").Value = "  -5.36%  "

$ws.Range("D41").Value = "'0.0831"
$ws.Range("E41").Value = "  -3.18%  "

$ws.Range("D42").Value = "'5.05"
$ws.Range("E42").Value = "  -3.50%  "

$ws.Range("D43").Value = "'0.870"
$ws.Range("E43").Value = "  -3.25%  "

$ws.Range("D44").Value = "'45.59"
$ws.Range("E44").Value = "  -1.52%  "

$ws.Range("D45").Value = "'1.80"
$ws.Range("E45").Value = "  -5.67%  "

$ws.Range("D46").Value = "'0.998"
$ws.Range("E46").Value = "  -0.14%  "

$ws.Range("E47").Value = "  -4.94%  "

$ws.Range("D48").Value = "'23.89"
$ws.Range("E48").Value = "  -1.63%  "

$ws.Range("D49").Value = "'7.17"
$ws.Range("E49").Value = "  -0.18%  "

$ws.Range("D50").Value = "'1.14"
$ws.Range("E50").Value = "  -4.67%  "

$ws.Range("D51").Value = "'0.926"
$ws.Range("E51").Value = "  -2.79%  "
